$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: conversation_12_07_2023__10_22_51
$ws.Range("B2").Value = 1689535950035
$ws.Range("C2").Value = 1689536250035

# Row 3: conversation_12_07_2023__09_53_17
$ws.Range("B3").Value = 1689536281159
$ws.Range("C3").Value = 1689536576890
$ws.Range("D3").Value = 4

# Row 4: conversation_11_07_2023__14_51_17
$ws.Range("B4").Value = 1689536578067
$ws.Range("C4").Value = 1689536878067
